$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.06603978468115329
$ws.Range("H2").Value = 98.34758886485763
$ws.Range("I2").Value = 25.81699312732994
$ws.Range("G3").Value = 0.05325878423939615
$ws.Range("H3").Value = 5.536714697772589
$ws.Range("G4").Value = 0.05107971399761535
$ws.Range("H4").Value = 4.938575083177073
$ws.Range("G5").Value = 0.05787459370931597
$ws.Range("H5").Value = -18.08793973580863
$ws.Range("G6").Value = -0.1193263051722239
$ws.Range("H6").Value = -1.445795045504012
$ws.Range("G7").Value = -0.1102156219180547
$ws.Range("H7").Value = 11.85110544247711
$ws.Range("G8").Value = -0.1813451552608865
$ws.Range("H8").Value = 8.953550851936246
$ws.Range("G9").Value = -0.2806981108175262
$ws.Range("H9").Value = 7.0772631988347
$ws.Range("G10").Value = -0.002691935776942437
$ws.Range("H10").Value = -311.5112731182062
$ws.Range("G11").Value = 0.04219137017150613
$ws.Range("H11").Value = 305.177168977795
$ws.Range("G12").Value = 0.2018591090105028
$ws.Range("H12").Value = -4.743574185682092
$ws.Range("G13").Value = 0.2348148175427455
$ws.Range("H13").Value = -0.07735084626877313
$ws.Range("G14").Value = -0.06270748657100526
$ws.Range("H14").Value = 31.12871929894201
$ws.Range("G15").Value = -0.06390788596983148
$ws.Range("H15").Value = 9.94923533354044
$ws.Range("G16").Value = 0.1668427290064805
$ws.Range("H16").Value = -12.85459828349269
$ws.Range("G17").Value = 0.2102598800818728
$ws.Range("H17").Value = 20.92962306658706
$ws.Range("G18").Value = 0.05849667852494526
$ws.Range("H18").Value = 7.868031677591109
$ws.Range("G19").Value = 0.06762543899312602
$ws.Range("H19").Value = -21.40196145076206
$ws.Range("G20").Value = -0.01670981480805889
$ws.Range("H20").Value = -231.2745382075581
$ws.Range("G21").Value = -0.04511704937241399
$ws.Range("H21").Value = 16.25417497302972
$ws.Range("G22").Value = 0.05601503543311501
$ws.Range("H22").Value = -14.1831101559336
$ws.Range("G23").Value = 0.08217775753201385
$ws.Range("H23").Value = 42.4881716335076
$ws.Range("G24").Value = 0.04890865883062724
$ws.Range("H24").Value = 50.98009611073638
$ws.Range("G25").Value = 0.02667230920103616
$ws.Range("H25").Value = -9.376899170292031
$ws.Range("G26").Value = 0.1083645803741487
$ws.Range("H26").Value = -4.348814926324891
$ws.Range("G27").Value = 0.09495499329867346
$ws.Range("H27").Value = 5.286018028822811
$ws.Range("G28").Value = 0.1347432607022712
$ws.Range("H28").Value = 14.70721705666124
$ws.Range("G29").Value = 0.1426034792118002
$ws.Range("H29").Value = 19.20610487258346
$ws.Range("G30").Value = 0.06642409821940501
$ws.Range("H30").Value = -1.19627440505757
$ws.Range("G31").Value = 0.04836861107598718
$ws.Range("H31").Value = -29.50996927018731
$ws.Range("G32").Value = 0.06329986925804701
$ws.Range("H32").Value = 44.96190017013478
$ws.Range("G33").Value = 0.07036845858977764
$ws.Range("H33").Value = 29.50040008213798
$ws.Range("G34").Value = -0.00144061226947688
$ws.Range("H34").Value = 92.45698176200813
$ws.Range("G35").Value = 0.05045094592216719
$ws.Range("H35").Value = 261.5863566729951
$ws.Range("G36").Value = -0.008634980055699945
$ws.Range("H36").Value = -155.8575986011582
$ws.Range("G37").Value = 0.02103459616907108
$ws.Range("H37").Value = 67.95870304335951
$ws.Range("G38").Value = 0.05108023211144321
$ws.Range("H38").Value = -28.80089385771608
$ws.Range("G39").Value = 0.06537890562396609
$ws.Range("H39").Value = 51.82327220928573
$ws.Range("G40").Value = 0.03990923495873493
$ws.Range("H40").Value = -10.77701030537473
$ws.Range("G41").Value = 0.05770758606525871
$ws.Range("H41").Value = 366.915458185029
$ws.Range("G42").Value = 0.07137511914022257
$ws.Range("H42").Value = 36.52737942947857
$ws.Range("G43").Value = 0.06945951351875042
$ws.Range("H43").Value = 39.19774187337769
$ws.Range("G44").Value = 0.08254840554808492
$ws.Range("H44").Value = -37.34507555347311
$ws.Range("G45").Value = 0.125015564099229
$ws.Range("H45").Value = -30.32498860469694
$ws.Range("G46").Value = -0.0185847096237442
$ws.Range("H46").Value = 57.69456759652364
$ws.Range("G47").Value = 0.01469345588740302
$ws.Range("H47").Value = 660.911672184669
$ws.Range("G48").Value = 0.009517001786506489
$ws.Range("H48").Value = -34.32846053033172
$ws.Range("G49").Value = -0.008092832778719163
$ws.Range("H49").Value = -45.58548493264819
$ws.Range("G50").Value = 0.1211209840834479
$ws.Range("H50").Value = -15.25940120046192
$ws.Range("G51").Value = 0.1583066076119766
$ws.Range("H51").Value = 20.87683596494371
$ws.Range("G52").Value = 0.08814188728498154
$ws.Range("H52").Value = 42.27286485722038
$ws.Range("G53").Value = 0.05663937652428736
$ws.Range("H53").Value = -7.424208857842618
$ws.Range("G54").Value = -0.1048322649155554
$ws.Range("H54").Value = -17.38423138286415
$ws.Range("G55").Value = -0.06706356776482356
$ws.Range("H55").Value = 35.38321803831052
$ws.Range("G56").Value = 0.1472916104119313
$ws.Range("H56").Value = -5.048885557279886
$ws.Range("G57").Value = 0.1745713732693138
$ws.Range("H57").Value = 25.15324207707652
